$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ihme_tab_1.2")

$rng = $ws.Range("G2:G21")
$rng.NumberFormat = "@"
$rng.Value = "24"
